$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B3").Value = "2.0.0"
$ws1.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$ws1.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Sheet "Include from FSIII" updates ---
$ws2 = $wb.Worksheets.Item("Include from FSIII")

# Extend formatting (style s="2") from the last existing data row (row 5)
# down through the two new rows (6 and 7) before touching any values.
$ws2.Range("A5:B5").Copy()
$ws2.Range("A6:B7").PasteSpecial(-4122)  # xlPasteFormats

# New bottom row: System URI / urn:oid:1.2.208.176.2.21
$ws2.Range("A7").Value = "System URI"
$ws2.Range("B7").Value = "urn:oid:1.2.208.176.2.21"

# Row 6 stays blank (both columns empty) - leave as pasted-format-only blank cells.

# Existing concept rows shift down in meaning: insert the two new UUID
# concepts right after the header, pushing "C" and "D" down one row each.
$ws2.Range("A5").Value = "D"
$ws2.Range("B5").Value = ""

$ws2.Range("A4").Value = "C"
$ws2.Range("B4").Value = ""

$ws2.Range("A3").Value = "687159ad-a61c-47c0-a878-53aa54bae2d5"
$ws2.Range("B3").Value = ""

$ws2.Range("A2").Value = "d6d48a71-b96f-4b88-86f9-b13bd3c03560"
$ws2.Range("B2").Value = ""
